$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit re-orders three observation records currently on rows 20-22
# (a re-sort that happened upstream): the data that was on row 21 moves to
# row 20, the data that was on row 22 moves to row 21, and the data that
# was on row 20 moves to row 22. We apply it as direct per-cell writes
# (rather than a row copy/paste) so each cell lands on the correct row.

function Set-TimeText($addr, $text) {
    # Clock-time-looking text (e.g. "16:06") gets reinterpreted by Excel as
    # a time serial number unless the cell is pre-formatted as Text.
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

function Add-EmptyCell($addr) {
    # Touch the cell (re-assign a property to its own value) so it becomes
    # a present-but-empty cell in the sheet instead of a wholly absent one,
    # without leaving it with any actual value or a new, visible style.
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $rng.Font.Bold
}

# ---------------------------------------------------------------
# Row 20 <- former row 21 content
# ---------------------------------------------------------------
$ws.Range("A20").Value = 111986518
$ws.Range("B20").Value = 90812
$ws.Range("D20").Value = "LC"
$ws.Range("E20").Value = 4366
$ws.Range("F20").Value = "Skarp dropptaggsvamp"
$ws.Range("G20").Value = "Hydnellum peckii"
$ws.Range("H20").Value = "Banker"
Add-EmptyCell "K20"
$ws.Range("P20").Value = "Mellandammen (Mellandammen), Dlr"
$ws.Range("Q20").Value = 396446
$ws.Range("R20").Value = 6849382
$ws.Range("S20").Value = 5
Set-TimeText "Z20" "16:06"
Set-TimeText "AB20" "16:06"
$ws.Range("AW20").Value = "Bo karlstens"
$ws.Range("AX20").Value = "Bo karlstens, Göran Ehn, Erik Danielsson, Kajsa Larsson, Lars-Erik Nilsson, Lisa Olson, Håkan Thenander"

# ---------------------------------------------------------------
# Row 21 <- former row 22 content
# ---------------------------------------------------------------
$ws.Range("A21").Value = 111986256
$ws.Range("B21").Value = 90794
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 4362
$ws.Range("F21").Value = "Blå taggsvamp"
$ws.Range("G21").Value = "Hydnellum caeruleum"
$ws.Range("H21").Value = "(Hornem.) P.Karst."
$ws.Range("Q21").Value = 396466
$ws.Range("R21").Value = 6849432
Set-TimeText "Z21" "15:50"
Set-TimeText "AB21" "15:50"
# S21, K21, P21, AW21, AX21 keep their current values (unchanged by the diff)

# ---------------------------------------------------------------
# Row 22 <- former row 20 content
# ---------------------------------------------------------------
$ws.Range("A22").Value = 112037911
$ws.Range("B22").Value = 78202
# D22 stays "NT" (unchanged by the diff)
$ws.Range("E22").Value = 229821
$ws.Range("F22").Value = "Vedflamlav"
$ws.Range("G22").Value = "Ramboldia elabens"
$ws.Range("H22").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("K22").ClearContents() | Out-Null
$ws.Range("P22").Value = "Mellandammen, Dlr"
$ws.Range("Q22").Value = 396463
$ws.Range("R22").Value = 6849431
$ws.Range("S22").Value = 1
$ws.Range("Z22").ClearContents() | Out-Null
$ws.Range("AB22").ClearContents() | Out-Null
$ws.Range("AW22").Value = "Lars-Erik Nilsson"
$ws.Range("AX22").Value = "Lars-Erik Nilsson, Bo karlstens, Erik Danielsson, Göran Ehn, Håkan Thenander, Kajsa Larsson, Lisa Olson"
